# Auto-generated Excel COM-interop script updating the cryptos list
# per the commit "Updated cryptos list on Wed Feb 28 17:19:17 UTC 2024 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cellRef, $text)
    $r = $ws.Range($cellRef)
    # Force text entry so numeric-looking strings (e.g. "1.00", "425.90")
    # are not reinterpreted as numbers, and restore the default style
    # afterwards so no incidental formatting change is introduced.
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "63.538.40"
Set-TextCell "E2" "  +11.16%  "

# Row 3
Set-TextCell "D3" "3.475.86"
Set-TextCell "E3" "  +7.19%  "

# Row 4
Set-TextCell "D4" "0.997"
Set-TextCell "E4" "  -0.35%  "

# Row 5
Set-TextCell "D5" "425.90"
Set-TextCell "E5" "  +7.75%  "

# Row 6
Set-TextCell "D6" "116.53"
Set-TextCell "E6" "  +8.56%  "

# Row 7
Set-TextCell "D7" "0.604"
Set-TextCell "E7" "  +6.18%  "

# Row 8
Set-TextCell "D8" "0.997"
Set-TextCell "E8" "  -0.26%  "

# Row 9
Set-TextCell "D9" "0.658"
Set-TextCell "E9" "  +6.75%  "

# Row 10
Set-TextCell "D10" "0.118"
Set-TextCell "E10" "  +21.72%  "

# Row 11
Set-TextCell "D11" "41.15"
Set-TextCell "E11" "  +5.65%  "

# Row 12
Set-TextCell "E12" "  +1.47%  "

# Row 13
Set-TextCell "D13" "3.989.04"
Set-TextCell "E13" "  +6.20%  "

# Row 14
Set-TextCell "D14" "8.69"
Set-TextCell "E14" "  +7.28%  "

# Row 15
Set-TextCell "D15" "20.24"
Set-TextCell "E15" "  +6.98%  "

# Row 16
Set-TextCell "D16" "3.572.98"
Set-TextCell "E16" "  +10.12%  "

# Row 17
Set-TextCell "D17" "1.07"
Set-TextCell "E17" "  +3.19%  "

# Row 18
Set-TextCell "D18" "62.852.34"
Set-TextCell "E18" "  +10.32%  "

# Row 19
Set-TextCell "D19" "11.12"
Set-TextCell "E19" "  +0.73%  "

# Row 20
Set-TextCell "D20" "0.0000122"
Set-TextCell "E20" "  +15.52%  "

# Row 21
Set-TextCell "D21" "3.45"
Set-TextCell "E21" "  +3.60%  "

# Row 22
Set-TextCell "D22" "13.35"
Set-TextCell "E22" "  +3.13%  "

# Row 23
Set-TextCell "D23" "312.25"
Set-TextCell "E23" "  +4.90%  "

# Row 24
Set-TextCell "D24" "77.60"
Set-TextCell "E24" "  +4.86%  "

# Row 25
Set-TextCell "D25" "3.35"
Set-TextCell "E25" "  +6.27%  "

# Row 26
Set-TextCell "D26" "30.29"
Set-TextCell "E26" "  +8.71%  "

# Row 27
Set-TextCell "D27" "4.52"
Set-TextCell "E27" "  +2.96%  "

# Row 28
Set-TextCell "D28" "8.06"
Set-TextCell "E28" "  +4.79%  "

# Row 29
Set-TextCell "B29" "RenderToken"
Set-TextCell "C29" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D29" "7.70"
Set-TextCell "E29" "  +6.47%  "

# Row 30
Set-TextCell "B30" "Kaspa"
Set-TextCell "C30" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D30" "0.179"
Set-TextCell "E30" "  +6.45%  "

# Row 31
Set-TextCell "E31" "  +7.20%  "

# Row 32
Set-TextCell "D32" "11.69"
Set-TextCell "E32" "  +6.78%  "

# Row 33
Set-TextCell "B33" "Toncoin"
Set-TextCell "C33" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell "D33" "2.55"
Set-TextCell "E33" "  +21.00%  "

# Row 34
Set-TextCell "B34" "Dai"
Set-TextCell "C34" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell "D34" "1.00"
Set-TextCell "E34" "  +0.03%  "

# Row 35
Set-TextCell "D35" "40.30"
Set-TextCell "E35" "  +7.78%  "

# Row 36
Set-TextCell "D36" "0.0516"
Set-TextCell "E36" "  +6.82%  "

# Row 37
Set-TextCell "D37" "52.71"
Set-TextCell "E37" "  +1.94%  "

# Row 38
Set-TextCell "D38" "3.15"
Set-TextCell "E38" "  +3.63%  "

# Row 39
Set-TextCell "B39" "LidoDAOToken"
Set-TextCell "C39" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D39" "3.49"
Set-TextCell "E39" "  -0.33%  "

# Row 40
Set-TextCell "B40" "FirstDigitalUSD"
Set-TextCell "C40" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell "D40" "0.994"
Set-TextCell "E40" "  -0.57%  "

# Row 41
Set-TextCell "B41" "EnergySwap"
Set-TextCell "C41" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D41" "28.00"
Set-TextCell "E41" "  +27.76%  "

# Row 42
Set-TextCell "B42" "Monero"
Set-TextCell "C42" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D42" "139.27"
Set-TextCell "E42" "  +3.94%  "

# Row 43
Set-TextCell "B43" "Stellar"
Set-TextCell "C43" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D43" "0.125"
Set-TextCell "E43" "  +4.50%  "

# Row 44
Set-TextCell "B44" "ARBITRUM"
Set-TextCell "C44" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell "D44" "1.97"
Set-TextCell "E44" "  +4.08%  "

# Row 45
Set-TextCell "B45" "TheGraph"
Set-TextCell "C45" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextCell "D45" "0.296"
Set-TextCell "E45" "  +5.04%  "

# Row 46
Set-TextCell "B46" "NEARProtocol"
Set-TextCell "C46" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell "D46" "4.08"
Set-TextCell "E46" "  +2.58%  "

# Row 47
Set-TextCell "B47" "Celestia"
Set-TextCell "C47" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextCell "D47" "17.18"
Set-TextCell "E47" "  +1.74%  "

# Row 48
Set-TextCell "B48" "WEMIXToken"
Set-TextCell "C48" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D48" "2.31"
Set-TextCell "E48" "  +10.87%  "

# Row 49
Set-TextCell "D49" "2.201.44"
Set-TextCell "E49" "  +2.97%  "

# Row 50
Set-TextCell "D50" "2.01"
Set-TextCell "E50" "  -1.03%  "

# Row 51
Set-TextCell "D51" "6.34"
Set-TextCell "E51" "  +7.23%  "

